$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be interpreted as text so values like
# "249.39" or "0.671" aren't auto-coerced to numbers by COM assignment -
# matches the original inlineStr storage for these cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.068.51"
$ws.Range("E2").Value = "  +1.50%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.059.96"
$ws.Range("E3").Value = "  -1.87%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "249.39"
$ws.Range("E5").Value = "  -1.28%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.671"
$ws.Range("E6").Value = "  +2.43%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.07%  "

# Row 8 - Solana
$ws.Range("D8").Value = "55.84"
$ws.Range("E8").Value = "  +17.84%  "

# Row 9 - now OKB (was Cardano)
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "60.92"
$ws.Range("E9").Value = "  +0.96%  "

# Row 10 - now Cardano (was OKB)
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.385"
$ws.Range("E10").Value = "  +1.55%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  +6.07%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +5.93%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "15.05"
$ws.Range("E13").Value = "  +3.31%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.350.04"
$ws.Range("E14").Value = "  -2.18%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.818"
$ws.Range("E15").Value = "  -1.54%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "5.30"
$ws.Range("E16").Value = "  +4.38%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.046.74"
$ws.Range("E17").Value = "  -2.07%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.007.28"
$ws.Range("E18").Value = "  +1.45%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0932"
$ws.Range("E19").Value = "  +12.43%  "

# Row 20 - Litecoin
$ws.Range("D20").Value = "73.51"
$ws.Range("E20").Value = "  +1.02%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "14.29"
$ws.Range("E21").Value = "  +8.53%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "5.36"
$ws.Range("E22").Value = "  +3.29%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "237.45"
$ws.Range("E23").Value = "  -1.03%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.09%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  -2.60%  "

# Row 26 - Monero
$ws.Range("D26").Value = "170.65"
$ws.Range("E26").Value = "  +0.13%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "9.03"
$ws.Range("E27").Value = "  -1.12%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "20.19"
$ws.Range("E28").Value = "  -5.73%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "1.98"
$ws.Range("E29").Value = "  -0.04%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  +1.48%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "4.63"
$ws.Range("E31").Value = "  +3.59%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "1.06"
$ws.Range("E32").Value = "  +6.84%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.0628"
$ws.Range("E33").Value = "  +2.02%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +7.18%  "

# Row 35 - Kaspa
$ws.Range("D35").Value = "0.0885"
$ws.Range("E35").Value = "  -2.50%  "

# Row 36 - BinanceUSD
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.13%  "

# Row 37 - LidoDAOToken
$ws.Range("D37").Value = "2.28"
$ws.Range("E37").Value = "  -6.25%  "

# Row 38 - WEMIXToken
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -4.89%  "

# Row 39 - TrustWalletToken
$ws.Range("E39").Value = "  +0.63%  "

# Row 40 - Cronos
$ws.Range("E40").Value = "  +22.47%  "

# Row 41 - InjectiveProtocol
$ws.Range("D41").Value = "17.75"
$ws.Range("E41").Value = "  +11.81%  "

# Row 42 - VeChain
$ws.Range("D42").Value = "0.0225"
$ws.Range("E42").Value = "  +0.93%  "

# Row 43 - ARBITRUM
$ws.Range("D43").Value = "1.15"
$ws.Range("E43").Value = "  -1.96%  "

# Row 44 - Aave
$ws.Range("D44").Value = "96.96"
$ws.Range("E44").Value = "  -0.59%  "

# Row 45 - HuobiToken
$ws.Range("E45").Value = "  +1.36%  "

# Row 46 - FTXToken
$ws.Range("D46").Value = "4.07"
$ws.Range("E46").Value = "  +46.87%  "

# Row 47 - Gas
$ws.Range("D47").Value = "13.82"
$ws.Range("E47").Value = "  -52.11%  "

# Row 48 - now RenderToken (was Maker)
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +7.04%  "

# Row 49 - now Maker (was RenderToken)
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.297.51"
$ws.Range("E49").Value = "  -2.23%  "

# Row 50 - MXToken
$ws.Range("D50").Value = "2.92"
$ws.Range("E50").Value = "  +2.56%  "

# Row 51 - THORChain
$ws.Range("D51").Value = "4.10"
$ws.Range("E51").Value = "  +7.37%  "

# Restore default (unstyled) formatting on the Price column now that the
# values are committed as text, so styling matches the original workbook.
$ws.Range("D2:D51").Style = "Normal"
